$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 281.7619
$ws.Range("I55").Value = 181.09091
$ws.Range("J55").Value = 392.5
$ws.Range("K55").Value = 181.09091
$ws.Range("L55").Value = 392.5
$ws.Range("M55").Value = 32.90908999999999
$ws.Range("N55").Value = -820.5
$ws.Range("H86").Value = 2281072
$ws.Range("I86").Value = 2496.6667
$ws.Range("K86").Value = 2496.6667
$ws.Range("M86").Value = -1373.6667
$ws.Range("H88").Value = 5051559
$ws.Range("I88").Value = 1003
$ws.Range("J88").Value = 5348650.5
$ws.Range("K88").Value = 1003
$ws.Range("L88").Value = 5348650.5
$ws.Range("M88").Value = -597
$ws.Range("N88").Value = -5349462.5
$ws.Range("H89").Value = 2281072
$ws.Range("I89").Value = 2496.6667
$ws.Range("K89").Value = 12483.3335
$ws.Range("M89").Value = -6867.333500000001
$ws.Range("H91").Value = 5051559
$ws.Range("I91").Value = 1003
$ws.Range("J91").Value = 5348650.5
$ws.Range("K91").Value = 1003
$ws.Range("L91").Value = 5348650.5
$ws.Range("M91").Value = 401
$ws.Range("N91").Value = -5351458.5
$ws.Range("H107").Value = 1349.5625
$ws.Range("H137").Value = 21443468
$ws.Range("I137").Value = 5682767
$ws.Range("J137").Value = 48115424
$ws.Range("K137").Value = 17048301
$ws.Range("L137").Value = 144346272
$ws.Range("M137").Value = -17045751
$ws.Range("N137").Value = -144351372
$ws.Range("H138").Value = 2257.5225
$ws.Range("I138").Value = 1603.5
$ws.Range("J138").Value = 3508.6956
$ws.Range("K138").Value = 4810.5
$ws.Range("L138").Value = 10526.0868
$ws.Range("M138").Value = 329.5
$ws.Range("N138").Value = -20806.0868
$ws.Range("H141").Value = 1383.1818
$ws.Range("I141").Value = 1341.5
$ws.Range("K141").Value = 4024.5
$ws.Range("M141").Value = 1155.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11820.5625
$ws.Range("I2").Value = 15503.25
$ws.Range("K2").Value = 15503.25
$ws.Range("M2").Value = -15390.25
$ws.Range("H45").Value = 257518.23
$ws.Range("I45").Value = 477040.75
$ws.Range("J45").Value = 1408.6111
$ws.Range("K45").Value = 477040.75
$ws.Range("L45").Value = 1408.6111
$ws.Range("M45").Value = -476663.75
$ws.Range("N45").Value = -2162.6111
$ws.Range("H97").Value = 633.3077
$ws.Range("I97").Value = 492.57144
$ws.Range("J97").Value = 797.5
$ws.Range("K97").Value = 492.57144
$ws.Range("L97").Value = 797.5
$ws.Range("M97").Value = 3.428560000000004
$ws.Range("N97").Value = -1789.5
$ws.Range("H116").Value = 11820.5625
$ws.Range("I116").Value = 15503.25
$ws.Range("K116").Value = 15503.25
$ws.Range("M116").Value = -13209.25
$ws.Range("H122").Value = 1835.5714
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 1933.1666
$ws.Range("K122").Value = 3750
$ws.Range("L122").Value = 5799.4998
$ws.Range("M122").Value = -1300
$ws.Range("N122").Value = -10699.4998
$ws.Range("H132").Value = 18524950
$ws.Range("I132").Value = 23817742
$ws.Range("J132").Value = 7939365
$ws.Range("K132").Value = 71453226
$ws.Range("L132").Value = 23818095
$ws.Range("M132").Value = -71450696
$ws.Range("N132").Value = -23823155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11820.5625
$ws.Range("I3").Value = 15503.25
$ws.Range("K3").Value = 15503.25
$ws.Range("M3").Value = -15389.25
$ws.Range("H86").Value = 1922.99
$ws.Range("I86").Value = 1935.1753
$ws.Range("J86").Value = 1529
$ws.Range("K86").Value = 1935.1753
$ws.Range("L86").Value = 1529
$ws.Range("M86").Value = -812.1753000000001
$ws.Range("N86").Value = -3775
$ws.Range("H89").Value = 1922.99
$ws.Range("I89").Value = 1935.1753
$ws.Range("J89").Value = 1529
$ws.Range("K89").Value = 9675.8765
$ws.Range("L89").Value = 7645
$ws.Range("M89").Value = -4059.8765
$ws.Range("N89").Value = -18877
$ws.Range("H94").Value = 1806.6538
$ws.Range("I94").Value = 1297
$ws.Range("J94").Value = 3190
$ws.Range("K94").Value = 1297
$ws.Range("L94").Value = 3190
$ws.Range("M94").Value = -846
$ws.Range("N94").Value = -4092
$ws.Range("H99").Value = 1707.3334
$ws.Range("I99").Value = 1496.5834
$ws.Range("J99").Value = 1988.3334
$ws.Range("K99").Value = 1496.5834
$ws.Range("L99").Value = 1988.3334
$ws.Range("M99").Value = 1.416600000000017
$ws.Range("N99").Value = -4984.3334
$ws.Range("H107").Value = 323461.97
$ws.Range("I107").Value = 769971.0600000001
$ws.Range("J107").Value = 983.1667
$ws.Range("K107").Value = 769971.0600000001
$ws.Range("L107").Value = 983.1667
$ws.Range("M107").Value = -768051.0600000001
$ws.Range("N107").Value = -4823.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3160024.2
$ws.Range("I31").Value = 1737850.2
$ws.Range("J31").Value = 6952488.5
$ws.Range("K31").Value = 1737850.2
$ws.Range("L31").Value = 6952488.5
$ws.Range("M31").Value = -1737555.2
$ws.Range("N31").Value = -6953078.5
$ws.Range("H34").Value = 3160024.2
$ws.Range("I34").Value = 1737850.2
$ws.Range("J34").Value = 6952488.5
$ws.Range("K34").Value = 1737850.2
$ws.Range("L34").Value = 6952488.5
$ws.Range("M34").Value = -1737648.2
$ws.Range("N34").Value = -6952892.5
$ws.Range("H62").Value = 2950
$ws.Range("I62").Value = 2487.5
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 2487.5
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -1863.5
$ws.Range("N62").Value = -6048
$ws.Range("H65").Value = 2950
$ws.Range("I65").Value = 2487.5
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 12437.5
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -9317.5
$ws.Range("N65").Value = -30240

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1006.04346
$ws.Range("J131").Value = 1111
$ws.Range("L131").Value = 3333
$ws.Range("N131").Value = -13413

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3734.5908
$ws.Range("I22").Value = 2841.111
$ws.Range("J22").Value = 4353.154
$ws.Range("K22").Value = 2841.111
$ws.Range("L22").Value = 4353.154
$ws.Range("M22").Value = -2546.111
$ws.Range("N22").Value = -4943.154
$ws.Range("H27").Value = 3734.5908
$ws.Range("I27").Value = 2841.111
$ws.Range("J27").Value = 4353.154
$ws.Range("K27").Value = 2841.111
$ws.Range("L27").Value = 4353.154
$ws.Range("M27").Value = -2734.111
$ws.Range("N27").Value = -4567.154
$ws.Range("H46").Value = 507.75
$ws.Range("I46").Value = 466.5
$ws.Range("J46").Value = 549
$ws.Range("K46").Value = 466.5
$ws.Range("L46").Value = 549
$ws.Range("M46").Value = -278.5
$ws.Range("N46").Value = -925
$ws.Range("H68").Value = 2927.2727
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 2970
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 2970
$ws.Range("M68").Value = -1751
$ws.Range("N68").Value = -4468
$ws.Range("H71").Value = 2927.2727
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 2970
$ws.Range("K71").Value = 12500
$ws.Range("L71").Value = 14850
$ws.Range("M71").Value = -8756
$ws.Range("N71").Value = -22338
$ws.Range("H82").Value = 4363.4
$ws.Range("I82").Value = 1523.5
$ws.Range("J82").Value = 5580.5
$ws.Range("K82").Value = 1523.5
$ws.Range("L82").Value = 5580.5
$ws.Range("M82").Value = -1162.5
$ws.Range("N82").Value = -6302.5
$ws.Range("H85").Value = 4363.4
$ws.Range("I85").Value = 1523.5
$ws.Range("J85").Value = 5580.5
$ws.Range("K85").Value = 1523.5
$ws.Range("L85").Value = 5580.5
$ws.Range("M85").Value = -275.5
$ws.Range("N85").Value = -8076.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1228
$ws.Range("I122").Value = 1084.3334
$ws.Range("J122").Value = 1551.25
$ws.Range("K122").Value = 3253.0002
$ws.Range("L122").Value = 4653.75
$ws.Range("M122").Value = -803.0001999999999
$ws.Range("N122").Value = -9553.75
$ws.Range("H126").Value = 31252438
$ws.Range("I126").Value = 50000420
$ws.Range("J126").Value = 5801.6665
$ws.Range("K126").Value = 150001260
$ws.Range("L126").Value = 17404.9995
$ws.Range("M126").Value = -149998790
$ws.Range("N126").Value = -22344.9995
$ws.Range("H135").Value = 51057.5
$ws.Range("J135").Value = 51057.5
$ws.Range("L135").Value = 51057.5
$ws.Range("N135").Value = -61197.5
